# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" on every sheet/cell
#   that shows it (Overview!E/F and the zh-cn / de-de Status column).
# - The Status column is narrower afterwards (it auto-sized to the new,
#   shorter text), so shrink those same columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Range("E4").Value = "In Translation"
$overview.Range("F4").Value = "In Translation"
$overview.Columns("E").ColumnWidth = 12.5
$overview.Columns("F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Range("C4").Value = "In Translation"
$zhcn.Columns("C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Range("C4").Value = "In Translation"
$dede.Columns("C").ColumnWidth = 12.5
